# Relatorio Analytics - fix typos / rewrap text per commit "Atualizando o relatório Analytics"
$d = $word.ActiveDocument

# 1) "A [Diwine] oferece" -> merge the misspelled-word run back into plain text
#    (drops the now-unneeded proofErr spell-check markers around "Diwine").
$d.Content.Find.Execute(
    "Diwine oferece uma solução prática e acessível para o ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Diwine oferece uma solução prática e acessível para o ", 2)

# 2) "open [source])" -> merge "open " with "source" so the spell-check markers
#    around "source" are cleared too.
$d.Content.Find.Execute(
    "open source", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "open source", 2)

# 3) "stand [by] ... 150 [uA]" -> merge the whole sentence into one run, clearing
#    the proofErr markers around "by" and "uA".
$d.Content.Find.Execute(
    "- Corrente: 200uA a 500mA, em stand by de 100uA a 150 uA", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "- Corrente: 200uA a 500mA, em stand by de 100uA a 150 uA", 2)

# 4) Rewrite the tail of the storage-simulation paragraph.
$d.Content.Find.Execute(
    "aceitáveis, portanto para a simulação utilizaremos como base 3 situações, sendo uma um cenário de eficiência e as outras 2 cenários que o processo será comprometido, para isso utilizaremos uma função para o primeiro cenário representando um queda de temperatura e umidade:  f(x)= x – 3 sendo x a temperatura e f(y) = y + 5 sendo y a umidade, para o segundo cenário: f(x)= x + 6 e f(y) = y – 7.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "aceitáveis, evitando a criação de fungos, a evaporação do vinho e a degradação das barricas no processo de maturação do vinho.",
    2)
